$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 12: was the "ProtonMass" row, becomes the "AdductLabels" row.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "AdductLabels"
$ws.Range("B12").Value = "proton"
$ws.Range("C12").Value = "Rarely"
$ws.Range("D12").Value = "Labels for the AdductMasses. Should be separated by a comma with no space (ex. proton,sodium)"

# ---------------------------------------------------------------------------
# Row 13 (new): "AdductMasses" row, holding the old ProtonMass numeric value.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "AdductMasses"
$ws.Range("C13").Value = "Rarely"
$ws.Range("D13").Value = "Masses for the Adducts. Should be separated by a comma with no space (ex. proton,sodium)"

# Clone the look of an existing "Rarely" data row (font / fill / alignment /
# full medium border) onto the new row's label/update/description cells,
# then strip the top+bottom edges so it reads as a continuation of row 12
# rather than a fully boxed row.
$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)

foreach ($addr in @("A13", "C13", "D13")) {
  $r = $ws.Range($addr)
  $r.Borders.Item(8).LineStyle = -4142
  $r.Borders.Item(9).LineStyle = -4142
}

# B13 holds the bare numeric value with no special formatting (matches the
# author's un-styled paste of the proton mass constant).
$ws.Range("B13").Value = 1.0072764700000001

$ws.Rows.Item(13).RowHeight = 17

# ---------------------------------------------------------------------------
# Selection, matching the saved cursor position in the source workbook.
# ---------------------------------------------------------------------------
$ws.Range("D12:D13").Select()
